$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed new shared strings in the same order the original commit introduced
#     them (Species "RX", Species "YFS", then the two new edf-pair strings)
#     so the rebuilt sharedStrings table lines up with the source workbook. ---
$ws.Range("B42").Value = "RX"
$ws.Range("B47").Value = "YFS"
$ws.Range("L43").Value = "3.99; 6.112"
$ws.Range("J44").Value = "27.109; 26.940"

# --- Fill in newly-computed "pheno thr temp" (M) / "geo thr temp" (K) values
#     for the AP and PK egg thr.pheno / thr.geo rows (GAMs re-run) ---
$ws.Range("M13").Value = 2.048
$ws.Range("K14").Value = 2.121
$ws.Range("M23").Value = 1.0488
$ws.Range("K24").Value = 1.0842000000000001

# --- Row 42: RX / egg / eg.base ---
$ws.Range("C42").Value = "egg"
$ws.Range("D42").Value = "eg.base"
$ws.Range("E42").Value = 0.53300000000000003
$ws.Range("F42").Value = 4865.8999999999996
$ws.Range("G42").Value = 9651.8799999999992
$ws.Range("H42").Value = 27.808
$ws.Range("I42").Value = 5.5229999999999997

# --- Row 43: RX / egg / thr.pheno ---
$ws.Range("B43").Value = "RX"
$ws.Range("C43").Value = "egg"
$ws.Range("D43").Value = "thr.pheno"
$ws.Range("E43").Value = 0.55000000000000004
$ws.Range("F43").Value = 4817.5
$ws.Range("G43").Value = 9547.268
$ws.Range("H43").Value = 27.765000000000001
$ws.Range("M43").Value = 2.0169999999999999

# --- Row 44: RX / egg / thr.geo ---
$ws.Range("B44").Value = "RX"
$ws.Range("C44").Value = "egg"
$ws.Range("D44").Value = "thr.geo"
$ws.Range("E44").Value = 0.63200000000000001
$ws.Range("F44").Value = 4580.7
$ws.Range("G44").Value = 9004.2510000000002
$ws.Range("I44").Value = 7.8869999999999996
$ws.Range("K44").Value = 2.2850000000000001

# --- Row 45: RX / egg / vc.pheno ---
$ws.Range("B45").Value = "RX"
$ws.Range("C45").Value = "egg"
$ws.Range("D45").Value = "vc.pheno"
$ws.Range("E45").Value = 0.54200000000000004
$ws.Range("F45").Value = 4847.5
$ws.Range("G45").Value = 9602.6669999999995
$ws.Range("H45").Value = 27.756
$ws.Range("I45").Value = 5.1840000000000002
$ws.Range("O45").Value = 6.2750000000000004

# --- Row 46: RX / egg / vc.geo ---
$ws.Range("B46").Value = "RX"
$ws.Range("C46").Value = "egg"
$ws.Range("D46").Value = "vc.geo"
$ws.Range("E46").Value = 0.58799999999999997
$ws.Range("F46").Value = 4726.5
$ws.Range("G46").Value = 9317.1419999999998
$ws.Range("H46").Value = 20.838999999999999
$ws.Range("I46").Value = 8.0250000000000004
$ws.Range("N46").Value = 21.994

# --- Row 47: YFS / larvae / lv.base ---
$ws.Range("C47").Value = "larvae"
$ws.Range("D47").Value = "lv.base"
$ws.Range("E47").Value = 0.79200000000000004
$ws.Range("F47").Value = 4002.2
$ws.Range("G47").Value = 7927.0129999999999
$ws.Range("H47").Value = 26.266999999999999
$ws.Range("I47").Value = 5.8179999999999996

# --- Row 48: YFS / larvae / lv.add.sal ---
$ws.Range("B48").Value = "YFS"
$ws.Range("C48").Value = "larvae"
$ws.Range("D48").Value = "lv.add.sal"
$ws.Range("E48").Value = 0.79700000000000004
$ws.Range("F48").Value = 3982.8
$ws.Range("G48").Value = 7872.8720000000003
$ws.Range("H48").Value = 26.303999999999998
$ws.Range("I48").Value = 5.83
$ws.Range("P48").Value = 7.016

# --- Row 49: YFS / larvae / lv.add.temp ---
$ws.Range("B49").Value = "YFS"
$ws.Range("C49").Value = "larvae"
$ws.Range("D49").Value = "lv.add.temp"
$ws.Range("E49").Value = 0.81799999999999995
$ws.Range("F49").Value = 3848
$ws.Range("G49").Value = 7604.3320000000003
$ws.Range("H49").Value = 25.678999999999998
$ws.Range("I49").Value = 5.1749999999999998
$ws.Range("Q49").Value = 8.423

# --- Row 50: YFS / larvae / lv.temp.sal (row previously had no Date, now does) ---
$ws.Range("A49").Copy()
$ws.Range("A50").PasteSpecial(-4122)
$ws.Range("A50").Value = 44543
$ws.Range("B50").Value = "YFS"
$ws.Range("C50").Value = "larvae"
$ws.Range("D50").Value = "lv.temp.sal"
$ws.Range("E50").Value = 0.82299999999999995
$ws.Range("F50").Value = 3827.2
$ws.Range("G50").Value = 7548.3649999999998
$ws.Range("H50").Value = 25.768999999999998
$ws.Range("I50").Value = 5.1959999999999997
$ws.Range("P50").Value = 6.8079999999999998
$ws.Range("Q50").Value = 8.4949999999999992

# --- Row 51: YFS / larvae / lv.2d (row previously had no Date, now does) ---
$ws.Range("A49").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("A51").Value = 44543
$ws.Range("B51").Value = "YFS"
$ws.Range("C51").Value = "larvae"
$ws.Range("D51").Value = "lv.2d"
$ws.Range("E51").Value = 0.83799999999999997
$ws.Range("F51").Value = 3751.9
$ws.Range("G51").Value = 7358.0079999999998
$ws.Range("H51").Value = 25.966000000000001
$ws.Range("I51").Value = 5.0709999999999997
$ws.Range("R51").Value = 25.713000000000001

# --- Restore view state: selection moved to K25 (scrolled further down the table) ---
$ws.Range("K25").Select()
